$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_MyAccount")

$ws.Range("L2").Value = 'VerifyTitle: City Marketplace'
$ws.Range("L3").Value = 'VerifyElement: null'
$ws.Range("L4").Value = 'Click: null'
$ws.Range("L5").Value = 'VerifyElement: null'
$ws.Range("L6").Value = 'SetText: Randomemailid'
$ws.Range("L7").Value = 'SetText: 123456'
$ws.Range("L8").Value = 'Click: null'
$ws.Range("L9").Value = 'VerifyText: Akash Sangal'
$ws.Range("L10").Value = 'Click: null'
$ws.Range("L11").Value = 'VerifyElement: null'
$ws.Range("L12").Value = 'Click: null'
$ws.Range("L13").Value = 'VerifyTitle: My Account'
$ws.Range("L14").Value = 'VerifyMyAccountMenu: Account Dashboard'
$ws.Range("L15").Value = 'VerifyText: Account Dashboard'
$ws.Range("L16").Value = 'VerifyText: Account Information'
$ws.Range("L17").Value = 'VerifyText: Contact Information'
$ws.Range("L18").Value = 'VerifyText: Akash Sangal'
$ws.Range("L19").Value = 'VerifyText: Randomemailid'
$ws.Range("L20").Value = 'VerifyText: Edit'
$ws.Range("L21").Value = 'VerifyText: Change Password'
$ws.Range("L22").Value = 'VerifyText: Address Book'
$ws.Range("L23").Value = 'VerifyText: Manage Addresses'
$ws.Range("L24").Value = 'VerifyText: Default Billing Address'
$ws.Range("L25").Value = 'VerifyText: Edit Address'
$ws.Range("L26").Value = 'VerifyText: Edit Address'
$ws.Range("L27").Value = 'VerifyMyAccountMenu: Address Book'
$ws.Range("L28").Value = 'SelectMyAccountTab: Address Book'
$ws.Range("L29").Value = 'VerifyText: Add New Address'
$ws.Range("L30").Value = 'VerifyElement: null'
$ws.Range("L31").Value = 'VerifyElement: null'
$ws.Range("L32").Value = 'VerifyElement: null'
$ws.Range("L33").Value = 'VerifyElement: null'
$ws.Range("L34").Value = 'VerifyElement: null'
$ws.Range("L35").Value = 'VerifyElement: null'
$ws.Range("L36").Value = 'VerifyElement: null'
$ws.Range("L37").Value = 'VerifyElement: null'
$ws.Range("L38").Value = 'VerifyElement: null'
$ws.Range("L39").Value = 'VerifyElement: null'
$ws.Range("L40").Value = 'VerifyMyAccountMenu: Account Information'
$ws.Range("L41").Value = 'SelectMyAccountTab: Account Information'
$ws.Range("L42").Value = 'VerifyText: Edit Account Information'
$ws.Range("L43").Value = 'VerifyElement: null'
$ws.Range("L44").Value = 'VerifyElement: null'
$ws.Range("L45").Value = 'VerifyElement: null'
$ws.Range("L46").Value = 'VerifyElement: null'
$ws.Range("L47").Value = 'VerifyElement: null'
$ws.Range("L48").Value = 'VerifyMyAccountMenu: My Saved Cards'
$ws.Range("L49").Value = 'SelectMyAccountTab: My Saved Cards'
$ws.Range("L50").Value = 'VerifyText: My Saved Cards'
$ws.Range("L51").Value = 'VerifyText: You do not have any saved cards yet.'
$ws.Range("L52").Value = 'VerifyText: Add a new saved card'
$ws.Range("L53").Value = 'VerifyText: You must first add a full billing address before you can add saved cards.'
$ws.Range("L54").Value = 'VerifyMyAccountMenu: My Wish List'
$ws.Range("L55").Value = 'SelectMyAccountTab: My Wish List'
$ws.Range("L56").Value = 'VerifyText: My Wish List'
$ws.Range("L57").Value = 'VerifyText: You have no items in your wish list.'
$ws.Range("L58").Value = 'SelectMyAccountTab: Account Dashboard'
$ws.Range("L59").Value = 'Click: null'
$ws.Range("L60").Value = 'SetText: demo'
$ws.Range("L61").Value = 'SetText: Sangal'
$ws.Range("L62").Value = 'Click: null'
$ws.Range("L63").Value = 'VerifyText: demo Sangal'
$ws.Range("L64").Value = 'VerifyText: demo Sangal'
$ws.Range("L65").Value = 'Click: null'
$ws.Range("L66").Value = 'SetText: Akash'
$ws.Range("L67").Value = 'SetText: Sangal'
$ws.Range("L68").Value = 'Click: null'
$ws.Range("L69").Value = 'VerifyText: Akash Sangal'
$ws.Range("L70").Value = 'Click: Change Password'
$ws.Range("L71").Value = 'SetText: 123456'
$ws.Range("L72").Value = 'SetText: 1234567'
$ws.Range("L73").Value = 'SetText: 1234567'
$ws.Range("L74").Value = 'Click: null'
$ws.Range("L75").Value = 'Click: null'
$ws.Range("L76").Value = 'Click: null'
$ws.Range("L77").Value = 'Wait: 6000'
$ws.Range("L78").Value = 'Click: null'
$ws.Range("L79").Value = 'VerifyElement: null'
$ws.Range("L80").Value = 'SetText: Randomemailid'
$ws.Range("L81").Value = 'SetText: 1234567'
$ws.Range("L82").Value = 'Click: null'
$ws.Range("L83").Value = 'VerifyText: Akash Sangal'
$ws.Range("L84").Value = 'Click: null'
$ws.Range("L85").Value = 'VerifyElement: null'
$ws.Range("L86").Value = 'Click: null'
$ws.Range("L87").Value = 'Click: Change Password'
$ws.Range("L88").Value = 'SetText: 1234567'
$ws.Range("L89").Value = 'SetText: 123456'
$ws.Range("L90").Value = 'SetText: 123456'
$ws.Range("L91").Value = 'Click: null'
